$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 379.795
$ws.Range("C3").Value = 411.945
$ws.Range("C4").Value = 428.075
$ws.Range("C5").Value = 443.495
$ws.Range("C6").Value = 458.595
